# fix some fight bug
# The "Article" table rows for LP+2 / MP+2 / PP+2 had their AddXx / Effect
# values mixed up. Re-align each row so that:
#   - MP+2 row now grants AddMp (column D) and the "blueflash" effect
#   - PP+2 row now grants AddPp (column E) and the "redflash" effect
#   - LP+2 row now grants AddLp (column C) and the "yellowflash" effect
# Row 4 (Id 58000001): LP+2 -> MP+2 (AddMp=2, Effect=blueflash)
# Row 5 (Id 58000002): MP+2 -> PP+2 (AddPp=2, Effect=redflash)
# Row 6 (Id 58000003): PP+2 -> LP+2 (AddLp=2, Effect=yellowflash)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("B4").Value = "MP+2"
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = "blueflash"

# --- Row 5 ---
$ws.Range("B5").Value = "PP+2"
$ws.Range("D5").Value = $null
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "redflash"

# --- Row 6 ---
$ws.Range("B6").Value = "LP+2"
$ws.Range("E6").Value = $null
$ws.Range("C6").Value = 2
$ws.Range("F6").Value = "yellowflash"

# Update the active selection to match the authored workbook (F4 instead of F5)
$ws.Range("F4").Select()
